$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recompute the Ccl4-Ccr3 LR-pair table to also include the "ECs" cluster
# alongside the existing "M2" cluster, expanding row 2 into four rows that
# cover every Sending/Target cluster combination of ECs and M2 (per Dr Hou's
# advice). Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E..T the associated expression/specificity metrics.
$data = @(
    @("ECs", "Ccl4", "Ccr3", "ECs", 1, 0.3333333333333333, 335.6132, 1006.8396, 0.80464917790985, 0.80464917790985, 1, 0.3333333333333333, 0.438062, 1.314186, 0.6074000808827777, 0.6074000808827777, 147.0193896184, 1323.1745065656, 0.4887439757447035, 0.4887439757447035),
    @("ECs", "Ccl4", "Ccr3", "M2", 1, 0.3333333333333333, 335.6132, 1006.8396, 0.80464917790985, 0.80464917790985, 3, 1, 0.2831463333333333, 0.8494390000000001, 0.3925999191172223, 0.3925999191172223, 95.02764699826668, 855.2488229844, 0.3159052021651465, 0.3159052021651465),
    @("M2", "Ccl4", "Ccr3", "ECs", 3, 1, 81.47937800000001, 244.438134, 0.19535082209015, 0.19535082209015, 1, 0.3333333333333333, 0.438062, 1.314186, 0.6074000808827777, 0.6074000808827777, 35.69301928543601, 321.237173568924, 0.1186561051380742, 0.1186561051380742),
    @("M2", "Ccl4", "Ccr3", "M2", 3, 1, 81.47937800000001, 244.438134, 0.19535082209015, 0.19535082209015, 3, 1, 0.2831463333333333, 0.8494390000000001, 0.3925999191172223, 0.3925999191172223, 23.07058712298067, 207.635284106826, 0.07669471695207576, 0.07669471695207576)
)

$nRows = 4
$nCols = 20
$arr = New-Object 'object[,]' $nRows,$nCols
for ($r = 0; $r -lt $nRows; $r++) {
    for ($c = 0; $c -lt $nCols; $c++) {
        $arr[$r,$c] = $data[$r][$c]
    }
}

$ws.Range("A2:T5").Value = $arr
